$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"
$ws.Columns.Item(1).ColumnWidth = 23.875
$ws.Columns.Item(2).ColumnWidth = 8.125
$ws.Range("D39").Select()
